$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto price/volume figures (Coinranking snapshot refresh).
# Each entry is a cell reference plus its new text value; values are
# written with a leading apostrophe so Excel stores them as text (matching
# the original inlineStr cells) instead of auto-converting to numbers.
$updates = @(
    @{ Cell = "D2"; Value = "30.127.64" },
    @{ Cell = "E2"; Value = "  +0.13%  " },
    @{ Cell = "D3"; Value = "1.906.17" },
    @{ Cell = "E3"; Value = "  -0.41%  " },
    @{ Cell = "D4"; Value = "0.9985" },
    @{ Cell = "E4"; Value = "  -0.07%  " },
    @{ Cell = "D5"; Value = "0.8418" },
    @{ Cell = "E5"; Value = "  +4.85%  " },
    @{ Cell = "D6"; Value = "242.33" },
    @{ Cell = "E6"; Value = "  -0.49%  " },
    @{ Cell = "D7"; Value = "0.9989" },
    @{ Cell = "D8"; Value = "0.3297" },
    @{ Cell = "E8"; Value = "  +3.11%  " },
    @{ Cell = "D9"; Value = "26.78" },
    @{ Cell = "E9"; Value = "  +1.17%  " },
    @{ Cell = "D10"; Value = "0.07089" },
    @{ Cell = "E10"; Value = "  +1.78%  " },
    @{ Cell = "D11"; Value = "0.08087" },
    @{ Cell = "E11"; Value = "  +0.95%  " },
    @{ Cell = "D12"; Value = "0.7670" },
    @{ Cell = "E12"; Value = "  +1.86%  " },
    @{ Cell = "D13"; Value = "1.906.59" },
    @{ Cell = "E13"; Value = "  -1.80%  " },
    @{ Cell = "D14"; Value = "5.288" },
    @{ Cell = "E14"; Value = "  +0.72%  " },
    @{ Cell = "D15"; Value = "92.82" },
    @{ Cell = "E15"; Value = "  -1.03%  " },
    @{ Cell = "D16"; Value = "30.127.90" },
    @{ Cell = "E16"; Value = "  +0.06%  " },
    @{ Cell = "D17"; Value = "14.18" },
    @{ Cell = "E17"; Value = "  +0.71%  " },
    @{ Cell = "D18"; Value = "5.896" },
    @{ Cell = "E18"; Value = "  -1.63%  " },
    @{ Cell = "D19"; Value = "245.53" },
    @{ Cell = "E19"; Value = "  -1.61%  " },
    @{ Cell = "E20"; Value = "  -0.67%  " },
    @{ Cell = "D21"; Value = "2.157.23" },
    @{ Cell = "E21"; Value = "  +0.08%  " },
    @{ Cell = "D22"; Value = "0.9998" },
    @{ Cell = "E22"; Value = "  +0.04%  " },
    @{ Cell = "D23"; Value = "0.9984" },
    @{ Cell = "E23"; Value = "  -0.11%  " },
    @{ Cell = "D24"; Value = "7.027" },
    @{ Cell = "E24"; Value = "  +1.15%  " },
    @{ Cell = "D25"; Value = "0.1794" },
    @{ Cell = "E25"; Value = "  +26.26%  " },
    @{ Cell = "D26"; Value = "9.319" },
    @{ Cell = "E26"; Value = "  -0.24%  " },
    @{ Cell = "D27"; Value = "165.56" },
    @{ Cell = "E27"; Value = "  -2.27%  " },
    @{ Cell = "D28"; Value = "19.01" },
    @{ Cell = "E28"; Value = "  +0.03%  " },
    @{ Cell = "D29"; Value = "2.109" },
    @{ Cell = "E29"; Value = "  +1.86%  " },
    @{ Cell = "D30"; Value = "1.367" },
    @{ Cell = "E30"; Value = "  -1.38%  " },
    @{ Cell = "D31"; Value = "1.517" },
    @{ Cell = "E31"; Value = "  -0.75%  " },
    @{ Cell = "D32"; Value = "0.05955" },
    @{ Cell = "E32"; Value = "  +8.39%  " },
    @{ Cell = "D33"; Value = "4.306" },
    @{ Cell = "E33"; Value = "  -1.23%  " },
    @{ Cell = "D34"; Value = "4.087" },
    @{ Cell = "D35"; Value = "1.275" },
    @{ Cell = "E35"; Value = "  +0.56%  " },
    @{ Cell = "D36"; Value = "0.7343" },
    @{ Cell = "D37"; Value = "2.711" },
    @{ Cell = "E37"; Value = "  -0.64%  " },
    @{ Cell = "D38"; Value = "0.01930" },
    @{ Cell = "E38"; Value = "  -0.27%  " },
    @{ Cell = "D39"; Value = "2.788" },
    @{ Cell = "D40"; Value = "0.4460" },
    @{ Cell = "E40"; Value = "  -0.16%  " },
    @{ Cell = "D41"; Value = "73.17" },
    @{ Cell = "E41"; Value = "  -0.26%  " },
    @{ Cell = "D42"; Value = "5.961" },
    @{ Cell = "E42"; Value = "  -4.13%  " },
    @{ Cell = "D43"; Value = "0.8576" },
    @{ Cell = "E43"; Value = "  +2.87%  " },
    @{ Cell = "D44"; Value = "1.912" },
    @{ Cell = "E44"; Value = "  -0.12%  " },
    @{ Cell = "D45"; Value = "0.9986" },
    @{ Cell = "D46"; Value = "102.27" },
    @{ Cell = "E46"; Value = "  +1.39%  " },
    @{ Cell = "D47"; Value = "7.595" },
    @{ Cell = "E47"; Value = "  -0.44%  " },
    @{ Cell = "D48"; Value = "9.843" },
    @{ Cell = "E48"; Value = "  -0.36%  " },
    @{ Cell = "D49"; Value = "1.005.89" },
    @{ Cell = "E49"; Value = "  +1.78%  " },
    @{ Cell = "D50"; Value = "2.058.10" },
    @{ Cell = "E50"; Value = "  -0.26%  " },
    @{ Cell = "D51"; Value = "1.524" },
    @{ Cell = "E51"; Value = "  +0.91%  " }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = "'" + $u.Value
    $ws.Range($u.Cell).ClearFormats()
}
